# Halted produtiction for now
# - Rename the sheet to the generic 'Sheet1'
# - Clear the stray 'Cell A1' test value and align A1's formatting with
#   the rest of the header/label cells (drop the dark-blue banner style)
# - Re-save the river-name labels, which picked up mojibake (UTF-8 bytes
#   re-decoded as Latin-1) somewhere upstream

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook used to be scoped to a specific date range; generalize the name.
$ws.Name = "Sheet1"

# A4 already carries the plain bold/bordered label style (no fill). Clone that
# format onto A1 so the old dark-blue "banner" font/fill fall out of use.
$ws.Range("A4").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# A1 only ever held a placeholder test string - clear it out.
$ws.Range("A1").Value = ""

# Re-apply the river names; upstream re-encoding turned the accented
# characters into mojibake (UTF-8 bytes mis-decoded as Latin-1/cp1252).
$ws.Range("A10").Value = "IguaÃ§u"
$ws.Range("A12").Value = "JacuÃ­"
$ws.Range("A15").Value = "ParanaÃ­ba"
$ws.Range("A16").Value = "ParanÃ¡"
$ws.Range("A17").Value = "ParaÃ­ba do Sul"
$ws.Range("A18").Value = "ParnaÃ­ba"
$ws.Range("A19").Value = "SÃ£o Francisco"
$ws.Range("A21").Value = "TietÃª"
